# Update "Oni" disability-prevalence sheet:
#  - retitle the report (row 1, merged A1:I1)
#  - rename + refresh the "family with disabilities" data row (row 4)
#  - insert a new "disabilities Persons" data row (row 5)
#  - drop the old "Note" row, keep the "Source" row (now row 6)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Title row (row 1) -- new wording, merged across A1:I1
# ---------------------------------------------------------------------
$ws.Range("A1:I1").UnMerge()
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Oni Municipality"
$ws.Range("A1:I1").Merge()
$ws.Range("A1:I1").HorizontalAlignment = -4108
$ws.Range("A1:I1").VerticalAlignment = -4108
$ws.Range("A1:I1").WrapText = $true
$ws.Range("A1").Font.Name = "Arial"
$ws.Range("A1").Font.Size = 11
$ws.Range("A1").Font.Bold = $true
$ws.Rows.Item(1).RowHeight = 51
$ws.Rows.Item(2).RowHeight = 14.5

# ---------------------------------------------------------------------
# 2) Row 4 -- relabel "Number of disability persons" ->
#    "family with disabilities Persons " and replace the data with the
#    new figures (no more confidential "..." markers). The old bottom
#    border (that used to separate it from "Source") is dropped -- only
#    the top border remains now that a new row sits below it.
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "family with disabilities Persons "
$ws.Range("A4:I4").Borders.Item(9).LineStyle = 0

$row4vals = @(243, 241, 239, 243, 227, 211, 188, 190)
$cols = @("B", "C", "D", "E", "F", "G", "H", "I")
for ($i = 0; $i -lt 8; $i++) {
    $cell = $ws.Range($cols[$i] + "4")
    $cell.Value = $row4vals[$i]
    $cell.NumberFormat = "#\ ##0"
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
    $cell.Borders.Item(9).LineStyle = 0
    $cell.Borders.Item(8).LineStyle = 0
}
$ws.Rows.Item(4).RowHeight = 24.75

# ---------------------------------------------------------------------
# 3) Insert a new row 5 for "disabilities Persons " with its own data
#    (this pushes the old Source/Note rows down by one). A bottom
#    border now closes off this row instead.
# ---------------------------------------------------------------------
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).RowHeight = 21

$ws.Range("A5").Value = "disabilities Persons "
$ws.Range("A5").Font.Name = "Arial"
$ws.Range("A5").Font.Size = 10
$ws.Range("A5").HorizontalAlignment = -4131
$ws.Range("A5").VerticalAlignment = -4108
$ws.Range("A5").WrapText = $true
$ws.Range("A5").Borders.Item(9).LineStyle = 1
$ws.Range("A5").Borders.Item(9).Weight = 2
$ws.Range("A5").Borders.Item(9).ColorIndex = 64

$row5vals = @(260, 257, 255, 262, 244, 229, 205, 207)
for ($i = 0; $i -lt 8; $i++) {
    $cell = $ws.Range($cols[$i] + "5")
    $cell.Value = $row5vals[$i]
    $cell.NumberFormat = "#\ ##0"
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
}
# Only the last data column (I5) keeps a visible bottom border, matching
# the source workbook's styling.
$ws.Range("I5").Borders.Item(9).LineStyle = 1
$ws.Range("I5").Borders.Item(9).Weight = 2
$ws.Range("I5").Borders.Item(9).ColorIndex = 64

# ---------------------------------------------------------------------
# 4) The old "Source" row is now row 6 (merge auto-shifted to A6:H6);
#    the old "Note" row is now row 7 -- delete it, it is no longer used.
# ---------------------------------------------------------------------
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).RowHeight = 27.75

# ---------------------------------------------------------------------
# 5) Final selection, matching the saved workbook state
# ---------------------------------------------------------------------
$ws.Range("A1:I1").Select()
